$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 needs to be the literal text "9234328533" stored as a shared string
# (not an auto-converted number). Typing a pure-digit value directly would
# be auto-coerced to a number by Excel; using a TEXT() helper formula and
# then Paste Special > Values keeps it as text without stamping a
# quote-prefixed cell style on the cell.
$helper = $ws.Cells.Item(50, 10)
$helper.Formula = "=TEXT(9234328533,""0"")"
$helper.Copy()
$ws.Range("A2").PasteSpecial(-4163)
$helper.ClearContents()
$excel.CutCopyMode = $false

$ws.Range("B2").Value = "FEB2025"
$ws.Range("C2").Value = 355.41

$ws.Range("C4").Select()
